$d = $word.ActiveDocument

# ---------------------------------------------------------------------------
# 1. Update the letter date: "September 19, 2025" -> "September 21, 2025".
#    We locate the paragraph containing the date via Find, then rewrite the
#    whole paragraph's Range.Text so the existing run (and its xml:space
#    attribute) is reused rather than discarded.
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.ClearFormatting()
$found = $r.Find.Execute("September 19, 2025")
if ($found) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $r.Start -and $p.Range.End -ge $r.End) {
            $p.Range.Text = "September 21, 2025"
            break
        }
    }
}

# ---------------------------------------------------------------------------
# 2. Split the sender address line "909 Story Road, San Jose CA 95122" (the
#    one in the main body, not the "PROPERTY ADDRESS" table) into two
#    separate paragraphs: "909 Story Road" and "San Jose, CA 95122".
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.ClearFormatting()
$found = $r.Find.Execute("909 Story Road, San Jose CA 95122")
if ($found) {
    $r.Text = "909 Story Road" + [char]13 + "San Jose, CA 95122"
}

# ---------------------------------------------------------------------------
# 3. Remove the now-superfluous empty "NoSpacing" paragraph that used to
#    immediately follow "Board of Directors".
# ---------------------------------------------------------------------------
$r = $d.Content
$r.Find.ClearFormatting()
$found = $r.Find.Execute("Board of Directors")
if ($found) {
    $count = $d.Paragraphs.Count
    for ($i = 1; $i -le $count; $i++) {
        $p = $d.Paragraphs.Item($i)
        if ($p.Range.Start -le $r.Start -and $p.Range.End -ge $r.End) {
            $next = $d.Paragraphs.Item($i + 1)
            $next.Range.Delete()
            break
        }
    }
}
